$d = $word.ActiveDocument

# 1) Convert the field-code hyperlink "Teamwork Crouse @ SoftUni" into a real
#    w:hyperlink pointing at the SoftUni trainings page.
$teamworkField = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "practical-teamwork-sept-2017") {
        $teamworkField = $f
    }
}
if ($teamworkField -ne $null) {
    $fieldRange = $teamworkField.Result.Duplicate
    $para = $fieldRange.Paragraphs(1)
    $paraRange = $para.Range.Duplicate
    $teamworkField.Delete()
    $d.Hyperlinks.Add($paraRange, "https://softuni.bg/trainings/1799/practical-teamwork-sept-2017", [Type]::Missing, [Type]::Missing, "Teamwork Crouse @ SoftUni") | Out-Null
}

# 2) Fill in the Trello board placeholder with the real board URL.
$d.Content.Find.Execute("board - ……", $true, $false, $false, $false, $false, $true, 1, $false, "board - https://trello.com/b/UVbqbPfA/auctionproject", 2) | Out-Null

# 3) Remove the "Atanas Galchov" team member bullet entirely.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Atanas Galchov*") {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}
